$wb = $excel.ActiveWorkbook

# Add Sheet2 right after Sheet1, and make it the active sheet
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Row 2: column headers
$ws2.Range("A2").Value = "Time limit"
$ws2.Range("B2").Value = "Chuffed"
$ws2.Range("C2").Value = "Kissat"
$ws2.Range("D2").Value = "OR-Tools"

# Time-limit values, rows 3-6
$ws2.Range("A4").Value = "1 minute"
$ws2.Range("A3").Value = "10 seconds"
$ws2.Range("A5").Value = "5 minutes"
$ws2.Range("A6").Value = "10 minutes"

# Title row, merged across A1:D1 and centered
$ws2.Range("A1").Value = "Number of instances solved"
$ws2.Range("A1:D1").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("A1:D1").Merge()

$ws2.Activate()
